$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "31.235.89"
$ws.Range("D2").Style = "Normal"

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.981.93"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +6.00%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9990"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.11%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.8004"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +70.30%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "252.50"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.60%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.9992"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.10%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3378"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +17.71%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "25.60"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +16.15%  "

# Row 10
$ws.Range("E10").Value = "  +7.31%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.8343"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +15.69%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08090"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +4.29%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.987.27"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +6.31%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "100.03"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +4.32%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.448"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +6.34%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "273.17"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.81%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "31.235.06"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +3.00%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "13.83"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +6.65%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007908"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +5.78%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "2.249.12"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +6.58%  "

# Row 21
$ws.Range("E21").Value = "  +9.44%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.001"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.04%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.9987"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.16%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.917"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +11.19%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.625"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +6.58%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "164.27"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.63%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.1481"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +54.46%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "19.73"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +5.82%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.172"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +16.24%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.562"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +6.05%  "

# Row 31
$ws.Range("E31").Value = "  +2.75%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.554"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +8.65%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.322"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +5.68%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.05140"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +7.05%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.208"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +8.18%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7538"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +9.59%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.779"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.28%  "

# Row 38
$ws.Range("B38").Value = "Frax"
$ws.Range("C38").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.9985"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.21%  "

# Row 39
$ws.Range("B39").Value = "VeChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01997"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +6.49%  "

# Row 40
$ws.Range("B40").Value = "MXToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.902"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +3.40%  "

# Row 41
$ws.Range("B41").Value = "FraxShare"
$ws.Range("C41").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.577"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +6.14%  "

# Row 42
$ws.Range("B42").Value = "Aave"
$ws.Range("C42").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "78.12"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +5.28%  "

# Row 43
$ws.Range("B43").Value = "TheSandbox"
$ws.Range("C43").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.4636"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +10.01%  "

# Row 44
$ws.Range("B44").Value = "RenderToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.048"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +6.21%  "

# Row 45
$ws.Range("B45").Value = "TrustWalletToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.8502"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.71%  "

# Row 46
$ws.Range("B46").Value = "Quant"
$ws.Range("C46").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "105.10"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +4.56%  "

# Row 47
$ws.Range("B47").Value = "PaxDollar"
$ws.Range("C47").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.9993"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.00%  "

# Row 48
$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.954"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +4.41%  "

# Row 49
$ws.Range("B49").Value = "Aptos"
$ws.Range("C49").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.457"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +7.56%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.4272"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +9.28%  "

# Row 51
$ws.Range("B51").Value = "Elrond"
$ws.Range("C51").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "36.37"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +3.25%  "
